$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Tundra / Virtus Pro score-label pairs
$ws.Range("A5").Value = "Tundra (2 — 0) -"
$ws.Range("C5").Value = "Virtus Pro (2 — 0) -"
$ws.Range("A6").Value = "Tundra (2 — 1) -"
$ws.Range("C6").Value = "Virtus Pro (2 — 1) -"

# Update the XXX / Sunway score-label pairs
$ws.Range("A8").Value = "XXX (2 — 0) -"
$ws.Range("C8").Value = "Sunway (2 — 0) -"
$ws.Range("A9").Value = "XXX )2 — 1) -"
$ws.Range("C9").Value = "Sunway (2 — 1) -"

# Leave the selection on the last edited cell, as in the authored workbook
$ws.Range("C9").Select()
